$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 47, shifting rows 47:58 down to 48:59.
$ws.Rows.Item(47).Insert()

# Match the date-formatted style used by column D in the other rows
$ws.Cells.Item(47, 4).NumberFormat = $ws.Cells.Item(48, 4).NumberFormat

# Fill in the new row 47 values
$ws.Cells.Item(47, 1).Value = 8
$ws.Cells.Item(47, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(47, 3).Value = "Coquimbo"
$ws.Cells.Item(47, 4).Value = 44798
$ws.Cells.Item(47, 5).Value = 4
$ws.Cells.Item(47, 6).Value = 100114007
$ws.Cells.Item(47, 7).Value = "Jengibre"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 400
$ws.Cells.Item(47, 11).Value = 14000
$ws.Cells.Item(47, 12).Value = 15000
$ws.Cells.Item(47, 13).Value = 14500
$ws.Cells.Item(47, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(47, 15).Value = "Perú"
$ws.Cells.Item(47, 16).Value = 1115
$ws.Cells.Item(47, 17).Value = 13
$ws.Cells.Item(47, 18).Value = "Hortaliza"
